$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain percentage (e.g. "86%") must be forced to
# text format first, otherwise Excel auto-converts them into a numeric percent
# value (0.86) formatted as a percentage, which changes both the stored value
# and the cell style.
$percentCells = @(
    "H2", "H3", "H8", "H10", "H14", "H15", "H19", "H20", "H22", "H31", "H38", "H46"
)
foreach ($cell in $percentCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply all data updates from the meteocat automatic refresh (2026-02-13 20:50).
$ws.Range('E2').Value = '2026-02-13 20:48:41'
$ws.Range('H2').Value = '86%'
$ws.Range('I2').Value = '2.6 mm'
$ws.Range('E3').Value = '2026-02-13 20:48:44'
$ws.Range('H3').Value = '88%'
$ws.Range('I3').Value = '6.5 mm'
$ws.Range('E4').Value = '2026-02-13 20:48:47'
$ws.Range('J4').Value = '994.0 hPa'
$ws.Range('E5').Value = '2026-02-13 20:48:49'
$ws.Range('I5').Value = '1.8 mm'
$ws.Range('E6').Value = '2026-02-13 20:48:52'
$ws.Range('J6').Value = '994.0 hPa'
$ws.Range('E7').Value = '2026-02-13 20:48:55'
$ws.Range('J7').Value = '994.4 hPa'
$ws.Range('E8').Value = '2026-02-13 20:48:57'
$ws.Range('H8').Value = '81%'
$ws.Range('J8').Value = '994.3 hPa'
$ws.Range('N8').Value = '6.6 °C 20:28 TU'
$ws.Range('E9').Value = '2026-02-13 20:49:00'
$ws.Range('E10').Value = '2026-02-13 20:49:03'
$ws.Range('H10').Value = '88%'
$ws.Range('I10').Value = '20.3 mm'
$ws.Range('E11').Value = '2026-02-13 20:49:05'
$ws.Range('E12').Value = '2026-02-13 20:49:08'
$ws.Range('E13').Value = '2026-02-13 20:49:10'
$ws.Range('E14').Value = '2026-02-13 20:49:13'
$ws.Range('H14').Value = '84%'
$ws.Range('L14').Value = '29.5 km/h - 292º 20:16 TU'
$ws.Range('E15').Value = '2026-02-13 20:49:15'
$ws.Range('H15').Value = '77%'
$ws.Range('E16').Value = '2026-02-13 20:49:18'
$ws.Range('I16').Value = '13.0 mm'
$ws.Range('E17').Value = '2026-02-13 20:49:21'
$ws.Range('O17').Value = '0.4 °C'
$ws.Range('E18').Value = '2026-02-13 20:49:23'
$ws.Range('J18').Value = '994.2 hPa'
$ws.Range('E19').Value = '2026-02-13 20:49:26'
$ws.Range('H19').Value = '90%'
$ws.Range('E20').Value = '2026-02-13 20:49:29'
$ws.Range('H20').Value = '94%'
$ws.Range('I20').Value = '23.7 mm'
$ws.Range('E21').Value = '2026-02-13 20:49:32'
$ws.Range('J21').Value = '997.2 hPa'
$ws.Range('E22').Value = '2026-02-13 20:49:34'
$ws.Range('H22').Value = '92%'
$ws.Range('L22').Value = '49.3 km/h - 325º 20:15 TU'
$ws.Range('E23').Value = '2026-02-13 20:49:37'
$ws.Range('I23').Value = '10.8 mm'
$ws.Range('E24').Value = '2026-02-13 20:49:40'
$ws.Range('J24').Value = '995.2 hPa'
$ws.Range('L24').Value = '52.9 km/h - 305º 20:27 TU'
$ws.Range('E25').Value = '2026-02-13 20:49:42'
$ws.Range('I25').Value = '9.2 mm'
$ws.Range('E26').Value = '2026-02-13 20:49:45'
$ws.Range('E27').Value = '2026-02-13 20:49:48'
$ws.Range('E28').Value = '2026-02-13 20:49:51'
$ws.Range('J28').Value = '994.5 hPa'
$ws.Range('E29').Value = '2026-02-13 20:49:53'
$ws.Range('I29').Value = '14.2 mm'
$ws.Range('O29').Value = '11.0 °C'
$ws.Range('E30').Value = '2026-02-13 20:49:56'
$ws.Range('J30').Value = '993.9 hPa'
$ws.Range('E31').Value = '2026-02-13 20:49:59'
$ws.Range('H31').Value = '74%'
$ws.Range('I31').Value = '4.8 mm'
$ws.Range('J31').Value = '992.9 hPa'
$ws.Range('E32').Value = '2026-02-13 20:50:02'
$ws.Range('L32').Value = '46.1 km/h - 294º 20:20 TU'
$ws.Range('E33').Value = '2026-02-13 20:50:04'
$ws.Range('J33').Value = '996.0 hPa'
$ws.Range('E34').Value = '2026-02-13 20:50:07'
$ws.Range('E35').Value = '2026-02-13 20:50:10'
$ws.Range('I35').Value = '8.2 mm'
$ws.Range('J35').Value = '995.2 hPa'
$ws.Range('E36').Value = '2026-02-13 20:50:12'
$ws.Range('J36').Value = '994.1 hPa'
$ws.Range('L36').Value = '54.4 km/h - 324º 20:15 TU'
$ws.Range('O36').Value = '10.7 °C'
$ws.Range('E37').Value = '2026-02-13 20:50:15'
$ws.Range('J37').Value = '996.0 hPa'
$ws.Range('E38').Value = '2026-02-13 20:50:18'
$ws.Range('H38').Value = '79%'
$ws.Range('L38').Value = '33.8 km/h - 272º 20:19 TU'
$ws.Range('E39').Value = '2026-02-13 20:50:20'
$ws.Range('I39').Value = '19.5 mm'
$ws.Range('E40').Value = '2026-02-13 20:50:23'
$ws.Range('J40').Value = '997.6 hPa'
$ws.Range('E41').Value = '2026-02-13 20:50:26'
$ws.Range('J41').Value = '994.6 hPa'
$ws.Range('E42').Value = '2026-02-13 20:50:29'
$ws.Range('O42').Value = '10.9 °C'
$ws.Range('E43').Value = '2026-02-13 20:50:31'
$ws.Range('E44').Value = '2026-02-13 20:50:34'
$ws.Range('I44').Value = '8.3 mm'
$ws.Range('E45').Value = '2026-02-13 20:50:37'
$ws.Range('I45').Value = '0.8 mm'
$ws.Range('O45').Value = '5.6 °C'
$ws.Range('E46').Value = '2026-02-13 20:50:39'
$ws.Range('H46').Value = '89%'
$ws.Range('J46').Value = '995.3 hPa'
$ws.Range('L46').Value = '29.9 km/h - 323º 20:24 TU'
$ws.Range('M46').Value = '11.6 °C 20:29 TU'
$ws.Range('O46').Value = '9.0 °C'
